# LoginTest_ID_01.xlsx update ("Test case1 xcel updated")
#
# Header row (row 1): lower-cases the two existing headers and adds a
# third column header.
# Data row (row 2): keeps the phone-number value in column A (still
# entered with a leading apostrophe so it stays text instead of being
# reinterpreted as a number), swaps the emailed "password" value in
# column B for a new one, and adds a new expected-color value in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "expected_color"

# --- Data row ---
# Leading apostrophe forces text storage (quote-prefixed), matching the
# existing formatting already used for this cell.
$ws.Range("A2").Value = "'9545351058"
$ws.Range("B2").Value = "amruta@12"
$ws.Range("C2").Value = "rgba(56, 88, 152, 1)"

# --- Column widths for the (now three) used columns ---
# (input values chosen so the engine's pixel-rounded stored width lands as
# close as possible to the authored 13.5703125 / 14.140625 / 31.28515625)
$ws.Columns.Item(1).ColumnWidth = 12.666
$ws.Columns.Item(2).ColumnWidth = 13.3335
$ws.Columns.Item(3).ColumnWidth = 30.5005

# --- Selection moves to the newly added C2 cell ---
$ws.Range("C2").Select() | Out-Null
